$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.544.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "'2.086.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'234.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").Value = "'58.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'0.0780"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "'2.392.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'14.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "'21.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'0.786"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'5.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "'2.079.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'37.485.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'6.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "'69.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'0.0₃0822"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "'226.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "'169.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "'8.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "'1.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.87%  "
$ws.Range("E29").Value = "  +4.34%  "
$ws.Range("D30").Value = "'19.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'4.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "'0.0618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").Value = "'4.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").Value = "'2.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "'3.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("D37").Value = "'1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  -5.51%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "'0.0961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "'1.485.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "'97.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("E46").Value = "  -11.00%  "
$ws.Range("D47").Value = "'1.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "'15.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'3.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'46.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.93%  "
